$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "66.945.31"
$ws.Range("E2").Value = "  -0.62%  "

# Row 3
$ws.Range("D3").Value = "3.513.20"
$ws.Range("E3").Value = "  +0.78%  "

# Row 4
$ws.Range("E4").Value = "  -0.15%  "

# Row 5
$ws.Range("D5").Value = "'608.57"
$ws.Range("E5").Value = "  +0.63%  "

# Row 6
$ws.Range("D6").Value = "'147.79"
$ws.Range("E6").Value = "  -1.97%  "

# Row 7
$ws.Range("D7").Value = "3.513.28"
$ws.Range("E7").Value = "  +0.93%  "

# Row 8
$ws.Range("E8").Value = "  +0.01%  "

# Row 9
$ws.Range("D9").Value = "'0.478"
$ws.Range("E9").Value = "  -1.43%  "

# Row 10
$ws.Range("D10").Value = "'0.142"
$ws.Range("E10").Value = "  -0.83%  "

# Row 11
$ws.Range("D11").Value = "'7.97"
$ws.Range("E11").Value = "  +5.29%  "

# Row 12
$ws.Range("D12").Value = "'0.423"
$ws.Range("E12").Value = "  -1.93%  "

# Row 13
$ws.Range("D13").Value = "'0.0000217"
$ws.Range("E13").Value = "  +0.42%  "

# Row 14
$ws.Range("D14").Value = "4.107.48"
$ws.Range("E14").Value = "  +0.67%  "

# Row 15
$ws.Range("D15").Value = "'31.87"
$ws.Range("E15").Value = "  -0.55%  "

# Row 16
$ws.Range("D16").Value = "3.512.05"
$ws.Range("E16").Value = "  +0.53%  "

# Row 17
$ws.Range("D17").Value = "66.982.10"
$ws.Range("E17").Value = "  -0.48%  "

# Row 18
$ws.Range("E18").Value = "  -0.44%  "

# Row 19
$ws.Range("D19").Value = "'10.66"
$ws.Range("E19").Value = "  +7.81%  "

# Row 20
$ws.Range("D20").Value = "'6.46"
$ws.Range("E20").Value = "  -0.70%  "

# Row 21
$ws.Range("D21").Value = "'15.31"
$ws.Range("E21").Value = "  -0.44%  "

# Row 22
$ws.Range("D22").Value = "'438.13"
$ws.Range("E22").Value = "  -1.62%  "

# Row 23
$ws.Range("D23").Value = "'0.608"
$ws.Range("E23").Value = "  -2.88%  "

# Row 24
$ws.Range("D24").Value = "'79.32"
$ws.Range("E24").Value = "  +0.96%  "

# Row 25
$ws.Range("D25").Value = "3.659.59"
$ws.Range("E25").Value = "  +0.87%  "

# Row 26
$ws.Range("E26").Value = "  -0.06%  "

# Row 27
$ws.Range("E27").Value = "  -3.61%  "

# Row 28
$ws.Range("D28").Value = "'9.76"
$ws.Range("E28").Value = "  -1.87%  "

# Row 29
$ws.Range("D29").Value = "'8.26"
$ws.Range("E29").Value = "  -4.57%  "

# Row 30
$ws.Range("E30").Value = "  +0.52%  "

# Row 31
$ws.Range("E31").Value = "  -3.79%  "

# Row 32
$ws.Range("E32").Value = "  -1.42%  "

# Row 33
$ws.Range("E33").Value = "  -0.06%  "

# Row 34
$ws.Range("D34").Value = "'25.48"
$ws.Range("E34").Value = "  -0.43%  "

# Row 35
$ws.Range("B35").Value = "NEARProtocol"
$ws.Range("C35").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D35").Value = "'5.96"
$ws.Range("E35").Value = "  -2.51%  "

# Row 36
$ws.Range("B36").Value = "ImmutableX"
$ws.Range("C36").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D36").Value = "'1.81"
$ws.Range("E36").Value = "  -2.52%  "

# Row 37
$ws.Range("B37").Value = "Aptos"
$ws.Range("C37").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D37").Value = "'8.02"
$ws.Range("E37").Value = "  +0.53%  "

# Row 38
$ws.Range("B38").Value = "USDe"
$ws.Range("C38").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D38").Value = "'1.00"
$ws.Range("E38").Value = "  +0.03%  "

# Row 39
$ws.Range("B39").Value = "FirstDigitalUSD"
$ws.Range("C39").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D39").Value = "'0.999"
$ws.Range("E39").Value = "  -0.12%  "

# Row 40
$ws.Range("B40").Value = "Monero"
$ws.Range("C40").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D40").Value = "'173.02"
$ws.Range("E40").Value = "  -2.79%  "

# Row 41
$ws.Range("B41").Value = "Hedera"
$ws.Range("C41").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D41").Value = "'0.0892"
$ws.Range("E41").Value = "  -0.16%  "

# Row 42
$ws.Range("B42").Value = "Filecoin"
$ws.Range("C42").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D42").Value = "'5.42"

# Row 43
$ws.Range("B43").Value = "Stacks"
$ws.Range("C43").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D43").Value = "'2.07"
$ws.Range("E43").Value = "  -9.80%  "

# Row 44
$ws.Range("B44").Value = "Mantle"
$ws.Range("C44").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D44").Value = "'0.895"
$ws.Range("E44").Value = "  +0.82%  "

# Row 45
$ws.Range("B45").Value = "OKB"
$ws.Range("C45").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D45").Value = "'46.16"
$ws.Range("E45").Value = "  -0.64%  "

# Row 46
$ws.Range("B46").Value = "InjectiveProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D46").Value = "'27.72"
$ws.Range("E46").Value = "  -7.66%  "

# Row 47
$ws.Range("B47").Value = "ONDO"
$ws.Range("C47").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D47").Value = "'1.28"
$ws.Range("E47").Value = "  -1.17%  "

# Row 48
$ws.Range("B48").Value = "Cosmos"
$ws.Range("C48").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D48").Value = "'7.46"
$ws.Range("E48").Value = "  -1.80%  "

# Row 49
$ws.Range("B49").Value = "dogwifhat"
$ws.Range("C49").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D49").Value = "'2.46"
$ws.Range("E49").Value = "  -3.11%  "

# Row 50
$ws.Range("B50").Value = "SuiNetwork"
$ws.Range("C50").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D50").Value = "'0.990"
$ws.Range("E50").Value = "  +0.78%  "

# Row 51
$ws.Range("B51").Value = "TheGraph"
$ws.Range("C51").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D51").Value = "'0.247"
$ws.Range("E51").Value = "  -1.92%  "
